$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.288.58"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.679.81"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5272"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2707"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07517"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.529"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.673.59"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5808"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008494"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "26.336.63"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.921"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.204"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.823"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1245"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06535"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("E29").Value = "  +4.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.330"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.600"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.590"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.659"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.033"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6238"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.405"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.747"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.450"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.70%  "
$ws.Range("D39").Value = "1.113.51"
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01622"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8758"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "1.830.93"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000114"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.187"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05276"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.095"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4294"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
